# Applies the "Add files via upload" edit to excelTest.xlsx:
#  - E5 becomes "7.5"
#  - B9, B10, B11 become "19:43:0"
#  - A new row 12 is appended with a full set of values
#
# Several of the new values look like numbers ("7.5", "3.0", "5.0", "6.0",
# "7.0", "8", "9.0") but must be stored as plain text (shared-string) cells,
# exactly like every other data cell on this sheet. Assigning such a literal
# straight to Range.Value would make Excel coerce it into a numeric cell, so
# for those we stage the literal as a text formula in a scratch cell, copy
# it, and paste-special just the (already-text) value into the destination
# -- this keeps the cell a plain shared-string cell with no number format /
# style changes, matching how the rest of the sheet is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "ZZ1"

function Set-TextLiteral($ws, $addr, $text) {
    $ws.Range($scratch).Formula = '="' + $text + '"'
    $ws.Range($scratch).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $ws.Range($scratch).ClearContents()
}

# --- existing-cell edits -----------------------------------------------

Set-TextLiteral $ws "E5" "7.5"

$ws.Range("B9").Value  = "19:43:0"
$ws.Range("B10").Value = "19:43:0"
$ws.Range("B11").Value = "19:43:0"

# --- new row 12 ----------------------------------------------------------

Set-TextLiteral $ws "A12" "31-12-1899"
Set-TextLiteral $ws "B12" "gfdgdf"
Set-TextLiteral $ws "C12" "3.0"
$ws.Range("D12").Value = "INDEFINIDO"
Set-TextLiteral $ws "E12" "ewrewrweurwg"
Set-TextLiteral $ws "F12" "6.0"
Set-TextLiteral $ws "G12" "7.0"
Set-TextLiteral $ws "H12" "8"
$ws.Range("I12").Value = "SIN_ASIGNAR"
$ws.Range("J12").Value = "No"
